# Update "implementasi graph dan tree-revisi.pptx"
#
# The diagram on slide 8 labels each city node with its shortest-path
# distance in parentheses, e.g. "Pangandaran (24570)". This commit
# corrects two of those distances:
#   - "Oval 13" (Pangandaran node):  (24570) -> (24590)
#   - "Oval 15" (Ciamis node):       (19860) -> (19880)
#
# Only the trailing " (NNNNN)" run inside each oval is touched; the
# oval's city-name run and the unrelated "... dengan jarak 24570" text
# in the slide's title textbox are left untouched.

$p = $ppt.ActivePresentation
$slide = $p.Slides.Item(8)

function Set-DistanceText {
    param(
        $Slide,
        [string]$ShapeName,
        [string]$OldFragment,
        [string]$NewFragment
    )

    for ($i = 1; $i -le $Slide.Shapes.Count; $i++) {
        $shape = $Slide.Shapes.Item($i)
        if ($shape.Name -ne $ShapeName) { continue }
        if (-not $shape.HasTextFrame) { continue }

        $textRange = $shape.TextFrame.TextRange
        $fullText = $textRange.Text
        $zeroBasedIndex = $fullText.IndexOf($OldFragment)
        if ($zeroBasedIndex -lt 0) { continue }

        # TextRange.Characters is 1-indexed.
        $fragment = $textRange.Characters($zeroBasedIndex + 1, $OldFragment.Length)
        $fragment.Text = $NewFragment
        return
    }
}

Set-DistanceText $slide "Oval 13" " (24570)" " (24590)"
Set-DistanceText $slide "Oval 15" " (19860)" " (19880)"
